$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the existing "Project 1" entry (row 9, D9) to "Data Cleaning Project"
$ws.Range("D9").Value = "Data Cleaning Project"

# Row 10: mark days 19-20 as the EDA Project
$ws.Range("C10").Value = "19-20"
$ws.Range("D10").Value = "EDA Project"

# Row 11: day 20, also EDA Project
$ws.Range("C11").Value = 20
$ws.Range("D11").Value = "EDA Project"

# Update the view: scroll so row 16 is at top, and select A2:E39 (active cell A2)
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("A2:E39").Select()
